# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamp cells produced by the handback
# report generation for the 3bb10003-... source file across the Overview,
# zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-09-04 17:11:30"

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-09-04 17:11:25"
$wsZhCn.Range("K2").Value = "2016-09-04 17:11:42"

# de-de: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-09-04 17:11:30"
$wsDeDe.Range("K2").Value = "2016-09-04 17:11:50"
